$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.063063333333333
$ws.Range("H2").Value = 6.189190000000001
$ws.Range("I2").Value = 0.1875199417503197
$ws.Range("J2").Value = 0.1875199417503197
$ws.Range("M2").Value = 13.713764
$ws.Range("N2").Value = 41.141292
$ws.Range("O2").Value = 0.0901423721847377
$ws.Range("P2").Value = 0.0901423721847377
$ws.Range("Q2").Value = 28.29236367038667
$ws.Range("R2").Value = 254.63127303348
$ws.Range("S2").Value = 0.01690349238131765
$ws.Range("T2").Value = 0.01690349238131765
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.063063333333333
$ws.Range("H3").Value = 6.189190000000001
$ws.Range("I3").Value = 0.1875199417503197
$ws.Range("J3").Value = 0.1875199417503197
$ws.Range("N3").Value = 84.55600199999999
$ws.Range("O3").Value = 0.1852659027513629
$ws.Range("P3").Value = 0.1852659027513629
$ws.Range("Q3").Value = 58.14812911315333
$ws.Range("R3").Value = 523.33316201838
$ws.Range("S3").Value = 0.03474105129225596
$ws.Range("T3").Value = 0.03474105129225596
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.063063333333333
$ws.Range("H4").Value = 6.189190000000001
$ws.Range("I4").Value = 0.1875199417503197
$ws.Range("J4").Value = 0.1875199417503197
$ws.Range("M4").Value = 21.07704566666666
$ws.Range("N4").Value = 63.23113699999999
$ws.Range("O4").Value = 0.1385421898057586
$ws.Range("P4").Value = 0.1385421898057586
$ws.Range("Q4").Value = 43.48328008989222
$ws.Range("R4").Value = 391.34952080903
$ws.Range("S4").Value = 0.02597942336233759
$ws.Range("T4").Value = 0.02597942336233759
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.063063333333333
$ws.Range("H5").Value = 6.189190000000001
$ws.Range("I5").Value = 0.1875199417503197
$ws.Range("J5").Value = 0.1875199417503197
$ws.Range("M5").Value = 89.15834833333334
$ws.Range("N5").Value = 267.475045
$ws.Range("O5").Value = 0.5860495352581409
$ws.Range("P5").Value = 0.5860495352581407
$ws.Range("Q5").Value = 183.9393193070611
$ws.Range("R5").Value = 1655.45387376355
$ws.Range("S5").Value = 0.1098959747144085
$ws.Range("T5").Value = 0.1098959747144085
$ws.Range("I6").Value = 0.5238509469163369
$ws.Range("J6").Value = 0.5238509469163369
$ws.Range("M6").Value = 13.713764
$ws.Range("N6").Value = 41.141292
$ws.Range("O6").Value = 0.0901423721847377
$ws.Range("P6").Value = 0.0901423721847377
$ws.Range("Q6").Value = 79.03682862149866
$ws.Range("R6").Value = 711.3314575934879
$ws.Range("S6").Value = 0.04722116702625971
$ws.Range("T6").Value = 0.04722116702625971
$ws.Range("I7").Value = 0.5238509469163369
$ws.Range("J7").Value = 0.5238509469163369
$ws.Range("N7").Value = 84.55600199999999
$ws.Range("O7").Value = 0.1852659027513629
$ws.Range("P7").Value = 0.1852659027513629
$ws.Range("S7").Value = 0.09705171858761144
$ws.Range("T7").Value = 0.09705171858761144
$ws.Range("I8").Value = 0.5238509469163369
$ws.Range("J8").Value = 0.5238509469163369
$ws.Range("M8").Value = 21.07704566666666
$ws.Range("N8").Value = 63.23113699999999
$ws.Range("O8").Value = 0.1385421898057586
$ws.Range("P8").Value = 0.1385421898057586
$ws.Range("Q8").Value = 121.4737869343409
$ws.Range("R8").Value = 1093.264082409068
$ws.Range("S8").Value = 0.07257545731760953
$ws.Range("T8").Value = 0.07257545731760953
$ws.Range("I9").Value = 0.5238509469163369
$ws.Range("J9").Value = 0.5238509469163369
$ws.Range("M9").Value = 89.15834833333334
$ws.Range("N9").Value = 267.475045
$ws.Range("O9").Value = 0.5860495352581409
$ws.Range("P9").Value = 0.5860495352581407
$ws.Range("Q9").Value = 513.8482109942644
$ws.Range("R9").Value = 4624.63389894838
$ws.Range("S9").Value = 0.3070026039848563
$ws.Range("T9").Value = 0.3070026039848562
$ws.Range("G10").Value = 2.101774
$ws.Range("H10").Value = 6.305322
$ws.Range("I10").Value = 0.1910385065181404
$ws.Range("J10").Value = 0.1910385065181404
$ws.Range("M10").Value = 13.713764
$ws.Range("N10").Value = 41.141292
$ws.Range("O10").Value = 0.0901423721847377
$ws.Range("P10").Value = 0.0901423721847377
$ws.Range("Q10").Value = 28.823232617336
$ws.Range("R10").Value = 259.409093556024
$ws.Range("S10").Value = 0.01722066415617465
$ws.Range("T10").Value = 0.01722066415617465
$ws.Range("G11").Value = 2.101774
$ws.Range("H11").Value = 6.305322
$ws.Range("I11").Value = 0.1910385065181404
$ws.Range("J11").Value = 0.1910385065181404
$ws.Range("N11").Value = 84.55600199999999
$ws.Range("O11").Value = 0.1852659027513629
$ws.Range("P11").Value = 0.1852659027513629
$ws.Range("Q11").Value = 59.239202182516
$ws.Range("R11").Value = 533.152819642644
$ws.Range("S11").Value = 0.0353929213703554
$ws.Range("T11").Value = 0.0353929213703554
$ws.Range("G12").Value = 2.101774
$ws.Range("H12").Value = 6.305322
$ws.Range("I12").Value = 0.1910385065181404
$ws.Range("J12").Value = 0.1910385065181404
$ws.Range("M12").Value = 21.07704566666666
$ws.Range("N12").Value = 63.23113699999999
$ws.Range("O12").Value = 0.1385421898057586
$ws.Range("P12").Value = 0.1385421898057586
$ws.Range("Q12").Value = 44.29918657901266
$ws.Range("R12").Value = 398.692679211114
$ws.Range("S12").Value = 0.02646689303024485
$ws.Range("T12").Value = 0.02646689303024485
$ws.Range("G13").Value = 2.101774
$ws.Range("H13").Value = 6.305322
$ws.Range("I13").Value = 0.1910385065181404
$ws.Range("J13").Value = 0.1910385065181404
$ws.Range("M13").Value = 89.15834833333334
$ws.Range("N13").Value = 267.475045
$ws.Range("O13").Value = 0.5860495352581409
$ws.Range("P13").Value = 0.5860495352581407
$ws.Range("Q13").Value = 187.3906984099434
$ws.Range("R13").Value = 1686.51628568949
$ws.Range("S13").Value = 0.1119580279613655
$ws.Range("T13").Value = 0.1119580279613654
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 1.073675666666667
$ws.Range("H14").Value = 3.221027
$ws.Range("I14").Value = 0.0975906048152031
$ws.Range("J14").Value = 0.09759060481520311
$ws.Range("M14").Value = 13.713764
$ws.Range("N14").Value = 41.141292
$ws.Range("O14").Value = 0.0901423721847377
$ws.Range("P14").Value = 0.0901423721847377
$ws.Range("Q14").Value = 14.72413470520933
$ws.Range("R14").Value = 132.517212346884
$ws.Range("S14").Value = 0.008797048620985692
$ws.Range("T14").Value = 0.008797048620985693
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 1.073675666666667
$ws.Range("H15").Value = 3.221027
$ws.Range("I15").Value = 0.0975906048152031
$ws.Range("J15").Value = 0.09759060481520311
$ws.Range("N15").Value = 84.55600199999999
$ws.Range("O15").Value = 0.1852659027513629
$ws.Range("P15").Value = 0.1852659027513629
$ws.Range("Q15").Value = 30.26190727267267
$ws.Range("R15").Value = 272.357165454054
$ws.Range("S15").Value = 0.0180802115011401
$ws.Range("T15").Value = 0.01808021150114011
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 1.073675666666667
$ws.Range("H16").Value = 3.221027
$ws.Range("I16").Value = 0.0975906048152031
$ws.Range("J16").Value = 0.09759060481520311
$ws.Range("M16").Value = 21.07704566666666
$ws.Range("N16").Value = 63.23113699999999
$ws.Range("O16").Value = 0.1385421898057586
$ws.Range("P16").Value = 0.1385421898057586
$ws.Range("Q16").Value = 22.62991105752211
$ws.Range("R16").Value = 203.669199517699
$ws.Range("S16").Value = 0.01352041609556665
$ws.Range("T16").Value = 0.01352041609556665
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 1.073675666666667
$ws.Range("H17").Value = 3.221027
$ws.Range("I17").Value = 0.0975906048152031
$ws.Range("J17").Value = 0.09759060481520311
$ws.Range("M17").Value = 89.15834833333334
$ws.Range("N17").Value = 267.475045
$ws.Range("O17").Value = 0.5860495352581409
$ws.Range("P17").Value = 0.5860495352581407
$ws.Range("Q17").Value = 95.72714908569056
$ws.Range("R17").Value = 861.5443417712152
$ws.Range("S17").Value = 0.05719292859751066
$ws.Range("T17").Value = 0.05719292859751066